# Append a new row (row 27) to each of the four worksheets, mirroring the
# last existing row (row 26) but with an updated timestamp in column A.

$wb = $excel.ActiveWorkbook

$timeValue = 45813.49222222222

$rowsData = @{
    1 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x78"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 376
        I = 15
    }
    2 = @{
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x8C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 396
        I = 14
    }
    3 = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6D"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 109
        I = 3
    }
    4 = @{
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6D"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 109
        I = 3
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i]

    $ws.Cells.Item(27, 1).Value = $timeValue
    $ws.Cells.Item(27, 1).NumberFormat = $ws.Cells.Item(26, 1).NumberFormat

    $ws.Cells.Item(27, 2).Value = $data.B
    $ws.Cells.Item(27, 3).Value = $data.C
    $ws.Cells.Item(27, 4).Value = $data.D
    $ws.Cells.Item(27, 5).Value = $data.E
    $ws.Cells.Item(27, 6).Value = $data.F
    $ws.Cells.Item(27, 7).Value = $data.G
    $ws.Cells.Item(27, 8).Value = $data.H
    $ws.Cells.Item(27, 9).Value = $data.I
}
